{"js": "// Replace the 25 \"###\u00d7#=\" problem strings in the practice-sheet table with\n// their updated values, per the commit's regenerated numbers.\nconst replacements = [\n  [\"218\u00d78=\", \"588\u00d72=\"],\n  [\"685\u00d75=\", \"965\u00d74=\"],\n  [\"434\u00d73=\", \"653\u00d76=\"],\n  [\"781\u00d78=\", \"414\u00d75=\"],\n  [\"876\u00d74=\", \"353\u00d78=\"],\n  [\"875\u00d72=\", \"472\u00d79=\"],\n  [\"878\u00d74=\", \"668\u00d76=\"],\n  [\"999\u00d77=\", \"824\u00d78=\"],\n  [\"663\u00d78=\", \"558\u00d72=\"],\n  [\"234\u00d79=\", \"539\u00d76=\"],\n  [\"994\u00d77=\", \"917\u00d76=\"],\n  [\"775\u00d76=\", \"917\u00d79=\"],\n  [\"935\u00d77=\", \"860\u00d78=\"],\n  [\"970\u00d76=\", \"602\u00d78=\"],\n  [\"290\u00d73=\", \"880\u00d78=\"],\n  [\"725\u00d76=\", \"241\u00d73=\"],\n  [\"691\u00d79=\", \"745\u00d78=\"],\n  [\"741\u00d74=\", \"101\u00d75=\"],\n  [\"457\u00d75=\", \"720\u00d74=\"],\n  [\"524\u00d75=\", \"755\u00d76=\"],\n  [\"339\u00d74=\", \"691\u00d76=\"],\n  [\"502\u00d74=\", \"509\u00d78=\"],\n  [\"379\u00d72=\", \"465\u00d77=\"],\n  [\"276\u00d77=\", \"342\u00d72=\"],\n  [\"844\u00d72=\", \"860\u00d79=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 \"###x#=\" problem strings in the practice-sheet table with\n# their updated values, per the commit's regenerated numbers.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"218\u00d78=\", \"588\u00d72=\"),\n  @(\"685\u00d75=\", \"965\u00d74=\"),\n  @(\"434\u00d73=\", \"653\u00d76=\"),\n  @(\"781\u00d78=\", \"414\u00d75=\"),\n  @(\"876\u00d74=\", \"353\u00d78=\"),\n  @(\"875\u00d72=\", \"472\u00d79=\"),\n  @(\"878\u00d74=\", \"668\u00d76=\"),\n  @(\"999\u00d77=\", \"824\u00d78=\"),\n  @(\"663\u00d78=\", \"558\u00d72=\"),\n  @(\"234\u00d79=\", \"539\u00d76=\"),\n  @(\"994\u00d77=\", \"917\u00d76=\"),\n  @(\"775\u00d76=\", \"917\u00d79=\"),\n  @(\"935\u00d77=\", \"860\u00d78=\"),\n  @(\"970\u00d76=\", \"602\u00d78=\"),\n  @(\"290\u00d73=\", \"880\u00d78=\"),\n  @(\"725\u00d76=\", \"241\u00d73=\"),\n  @(\"691\u00d79=\", \"745\u00d78=\"),\n  @(\"741\u00d74=\", \"101\u00d75=\"),\n  @(\"457\u00d75=\", \"720\u00d74=\"),\n  @(\"524\u00d75=\", \"755\u00d76=\"),\n  @(\"339\u00d74=\", \"691\u00d76=\"),\n  @(\"502\u00d74=\", \"509\u00d78=\"),\n  @(\"379\u00d72=\", \"465\u00d77=\"),\n  @(\"276\u00d77=\", \"342\u00d72=\"),\n  @(\"844\u00d72=\", \"860\u00d79=\")\n)\n\nforeach ($pair in $replacements) {\n  $old = $pair[0]\n  $new = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
